$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value/text updates (UC campus naming + merged record ids) ---
$ws.Range("A13").Value = 18070
$ws.Range("B13").Value = "University of California _ Berkeley"

$ws.Range("A15").Value = 23542
$ws.Range("B15").Value = "University of California _ Los Angeles"

$ws.Range("B74").Value = "University of California _ Davis"

$ws.Range("B80").Value = "University of California _ San Diego"

$ws.Range("B96").Value = "University of California _ Santa Barbara"

# --- Merging Scopus and masterlist: dedupe the two "Centre For/for Economic
#     Policy Research" rows (old rows 110 and 113) and insert a single
#     consolidated row (id 15030) ahead of the "University Of Massachusetts"
#     block, shifting the remaining rows up by one (old row 136 disappears).

# Delete the first duplicate (old row 110: 9778 / "Centre For Economic Policy Research")
$ws.Rows(110).Delete()

# The second duplicate (old row 113: 21413 / "Centre for Economic Policy Research ")
# is now at row 112 after the shift above - delete it too.
$ws.Rows(112).Delete()

# Insert the consolidated row ahead of row 106 ("University Of Massachusetts").
$ws.Rows(106).Insert()
$ws.Range("A106").Value = 15030
$ws.Range("B106").Value = "Centre For Economic Policy Research"

# The insert leaves A106 with a formatting gap (no border) - restore the
# same format used by the rest of the A column (id cells) from a neighbor.
$ws.Range("A107").Copy()
$ws.Range("A106").PasteSpecial(-4122)
